# Refresh the cryptos price/volume snapshot (and the row-38/39 coin swap)
# produced by the scheduled GitHub Actions scrape, cell-by-cell per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.048.31'
$ws.Range("E2").Value = '  -1.69%  '
# Row 3
$ws.Range("D3").Value = '1.767.88'
$ws.Range("E3").Value = '  -3.56%  '
# Row 4
$ws.Range("E4").Value = '  +0.28%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.83'
$ws.Range("E5").Value = '  -2.75%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9992'
$ws.Range("E6").Value = '  +0.57%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4278'
$ws.Range("E7").Value = '  -4.26%  '
# Row 8
$ws.Range("E8").Value = '  -4.47%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.59'
$ws.Range("E9").Value = '  -4.30%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07457'
$ws.Range("E10").Value = '  -4.28%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.094'
$ws.Range("E11").Value = '  -4.41%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9993'
$ws.Range("E12").Value = '  +0.45%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.18'
$ws.Range("E13").Value = '  -5.33%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.078'
$ws.Range("E14").Value = '  -4.47%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.337'
$ws.Range("E15").Value = '  -2.83%  '
# Row 16
$ws.Range("D16").Value = '1.788.69'
$ws.Range("E16").Value = '  -2.43%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.42'
$ws.Range("E17").Value = '  -1.16%  '
# Row 18
$ws.Range("E18").Value = '  -2.82%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06393'
$ws.Range("E19").Value = '  +0.39%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9998'
$ws.Range("E20").Value = '  +0.42%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.11'
$ws.Range("E21").Value = '  -3.17%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.973'
$ws.Range("E22").Value = '  -6.58%  '
# Row 23
$ws.Range("D23").Value = '28.042.55'
$ws.Range("E23").Value = '  -1.90%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.28'
$ws.Range("E24").Value = '  -5.15%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.131'
$ws.Range("E25").Value = '  -5.29%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.33'
$ws.Range("E26").Value = '  +2.14%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.22'
$ws.Range("E27").Value = '  -3.65%  '
# Row 28
$ws.Range("D28").Value = '1.980.86'
$ws.Range("E28").Value = '  -3.01%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.147'
$ws.Range("E29").Value = '  -10.46%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.15'
$ws.Range("E30").Value = '  -3.87%  '
# Row 31
$ws.Range("E31").Value = '  -5.40%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.648'
$ws.Range("E32").Value = '  -4.20%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.675'
$ws.Range("E33").Value = '  +0.04%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08921'
$ws.Range("E34").Value = '  -4.02%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.55'
$ws.Range("E35").Value = '  -2.70%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02314'
$ws.Range("E36").Value = '  -2.37%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2112'
$ws.Range("E37").Value = '  -4.32%  '
# Row 38
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06047'
$ws.Range("E38").Value = '  -3.73%  '
# Row 39
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.017'
$ws.Range("E39").Value = '  -3.93%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6390'
$ws.Range("E40").Value = '  -4.25%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.185'
$ws.Range("E41").Value = '  -1.06%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9986'
$ws.Range("E42").Value = '  +0.49%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.401'
$ws.Range("E43").Value = '  +0.01%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.837'
$ws.Range("E44").Value = '  -3.66%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.44'
$ws.Range("E45").Value = '  -4.59%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5932'
$ws.Range("E46").Value = '  -3.27%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.685'
$ws.Range("E47").Value = '  -2.43%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.008'
$ws.Range("E48").Value = '  -2.04%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '121.85'
$ws.Range("E49").Value = '  -4.72%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.181'
$ws.Range("E50").Value = '  +2.24%  '
# Row 51
$ws.Range("E51").Value = '  -2.31%  '
